# Weekly refresh of "Fruta / hortaliza" data: the data rows (2-48) get
# re-shuffled to a new row order. Build the new order by reading the
# current A2:T48 block and rewriting it back in permuted order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:T48")
$data = $rng.Value2

$nRows = $data.GetLength(0)
$nCols = $data.GetLength(1)

# For each new data-array row (1-based, row 1 == worksheet row 2), the
# index (1-based) of the data-array row whose values should be placed
# there. This mapping was derived from the target diff and is a full
# permutation of rows 2..48.
$srcIndex = @(3,4,19,21,22,34,35,39,40,12,42,13,18,7,25,26,23,44,45,16,38,36,1,29,14,15,43,6,28,27,37,30,2,8,9,10,5,11,31,32,33,46,47,24,41,17,20)

$newData = New-Object 'object[,]' $nRows,$nCols

for ($r = 1; $r -le $nRows; $r++) {
    $src = $srcIndex[$r - 1]
    for ($c = 1; $c -le $nCols; $c++) {
        $newData[$r-1, $c-1] = $data[$src, $c]
    }
}

$rng.Value2 = $newData
